# Add coal ramping parameters to the scenario settings sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells for the fossil-share columns (J1:K1) ---
$ws.Range("J1").Value = "MinFossilShare"
$ws.Range("K1").Value = "MaxFossilShare"

# --- "Assumed" / source note row (row 5) for the coal ramping columns ---
$ws.Range("L5").Value = "Assumed"
$ws.Range("M5").Value = "Assumed"
$ws.Range("N5").Value = "Source: A. Teruel, Perspestective of the Energy Transition: Technology Development and Investments under Uncertainty, Master thesis with DLR"

# --- Unit row (row 4) for the coal ramping columns ---
$ws.Range("L4").Value = "euros/MW"
$ws.Range("M4").Value = "euros/MW"
$ws.Range("N4").Value = "euros/MW"

# --- Header row (row 1) for the coal ramping columns ---
$ws.Range("L1").Value = "CoalRampingHourly"
$ws.Range("M1").Value = "CoalRampingDaily"
$ws.Range("N1").Value = "CoalRampingWearTear"

# --- Data row (row 2) values ---
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1
$ws.Range("J2:K2").NumberFormat = "0.0"

$ws.Range("L2").Value = 0.1
$ws.Range("M2").Value = 0.1
$ws.Range("N2").Value = 3.3

$ws.Range("K8").Select()
